# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns AD, AE, AF are appended after the existing "Unnamed: 28" (AC) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reuse the existing header formatting (bold, centered,
# thin border) by copying the format from the last header cell (AC1) onto the
# three new header cells, then set their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2 through 47): every player on the roster gets the team's
# overall record for the season: 86 wins, 77 losses, 0 ties.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 86   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 77   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
